$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "       Background del mapa de Haohamru(con todas sus animaciones)"
$ws.Range("C16").Value = "60 minuts"
$ws.Range("D16").Value = 110

$ws.Range("B25").Value = "Crear modulos menu y VictoryHaohmaru, y backgrounds"
$ws.Range("C25").Value = "45 minuts"
$ws.Range("D25").Value = 70

$ws.Range("E23").Select() | Out-Null
